# Generate Report for Handback
# Refresh the handoff/handback timestamps for the
# "9c788acf-9ee0-447c-a4ae-d1a226a7d50f" file row across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 9c788acf row (row 3)
$wsOverview.Range("G3").Value = "2016-08-19 06:46:42"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 9c788acf row (row 3)
$wsZhCn.Range("H3").Value = "2016-08-19 06:46:37"
$wsZhCn.Range("K3").Value = "2016-08-19 06:46:57"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 9c788acf row (row 3)
$wsDeDe.Range("H3").Value = "2016-08-19 06:46:42"
$wsDeDe.Range("K3").Value = "2016-08-19 06:47:12"
